# Investing Quotes.xlsx - "More work to Uranium ETF overview"
#
# Adds a new Charlie Munger quote row, and clears the stray/unused
# cell-level formatting that had accumulated on some of the quote rows
# (those cells carried an xf that applied no visible formatting; Excel
# collapses them back to the default style on the next save).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the vestigial per-cell styles (fontId=0/applyFont, applyAlignment
# with no real alignment) from the Author/Quote cells that had them.
$ws.Range("B7:B9").ClearFormats()
$ws.Range("C10:C14").ClearFormats()

# New quote: Charlie Munger, row 15.
$ws.Range("B15").Value = "Charlie Munger"
$ws.Range("C15").Value = "No wise pilot, no matter how great his talent and experience, fails to use his checklist"

# Leave the selection where the author last clicked after typing the entry.
$ws.Range("C16").Select() | Out-Null
